$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18
$ws.Range("O18").Value = 17
$ws.Range("P18").Value = 100
$ws.Range("Q18").Value = 3
$ws.Range("R18").Value = 12
$ws.Range("S18").Value = 1
$ws.Range("T18").Value = 0
$ws.Range("U18").Value = 2
$ws.Range("V18").Value = 0
$ws.Range("W18").Value = 1
$ws.Range("X18").Value = 100
$ws.Range("Y18").Value = $false
$ws.Range("Z18").Value = 0.1
$ws.Range("AA18").Value = 6
$ws.Range("AB18").Value = 200
$ws.Range("AC18").Value = 3
$ws.Range("AD18").Value = 7
$ws.Range("AE18").Value = 3000000
$ws.Range("AF18").Value = 0.16502394806634199

# Row 19
$ws.Range("O19").Value = 18
$ws.Range("P19").Value = 100
$ws.Range("Q19").Value = 3
$ws.Range("R19").Value = 12
$ws.Range("S19").Value = 1
$ws.Range("T19").Value = 1
$ws.Range("W19").Value = 1
$ws.Range("X19").Value = 100
$ws.Range("Z19").Value = 0.1
$ws.Range("AA19").Value = 6
$ws.Range("AB19").Value = 200
$ws.Range("AE19").Value = 3000000

# Row 20
$ws.Range("O20").Value = 19
$ws.Range("P20").Value = 100
$ws.Range("Q20").Value = 3
$ws.Range("R20").Value = 12
$ws.Range("S20").Value = 1
$ws.Range("T20").Value = 1
$ws.Range("W20").Value = 1
$ws.Range("X20").Value = 100
$ws.Range("Z20").Value = 0.1
$ws.Range("AA20").Value = 6
$ws.Range("AB20").Value = 200
$ws.Range("AE20").Value = 3000000

# Row 21
$ws.Range("O21").Value = 20
$ws.Range("P21").Value = 100
$ws.Range("Q21").Value = 3
$ws.Range("R21").Value = 12
$ws.Range("S21").Value = 1
$ws.Range("T21").Value = 1
$ws.Range("W21").Value = 1
$ws.Range("X21").Value = 100
$ws.Range("Z21").Value = 0.1
$ws.Range("AA21").Value = 6
$ws.Range("AB21").Value = 200
$ws.Range("AE21").Value = 3000000

# Row 22
$ws.Range("O22").Value = 21
$ws.Range("P22").Value = 100
$ws.Range("Q22").Value = 3
$ws.Range("R22").Value = 12
$ws.Range("S22").Value = 1
$ws.Range("T22").Value = 1
$ws.Range("W22").Value = 1
$ws.Range("X22").Value = 100
$ws.Range("Z22").Value = 0.1
$ws.Range("AA22").Value = 6
$ws.Range("AB22").Value = 200
$ws.Range("AE22").Value = 3000000

# Row 23
$ws.Range("O23").Value = 22
$ws.Range("P23").Value = 100
$ws.Range("Q23").Value = 3
$ws.Range("R23").Value = 12
$ws.Range("S23").Value = 1
$ws.Range("T23").Value = 1
$ws.Range("W23").Value = 1
$ws.Range("X23").Value = 100
$ws.Range("Z23").Value = 0.1
$ws.Range("AA23").Value = 6
$ws.Range("AB23").Value = 200
$ws.Range("AE23").Value = 3000000

# Row 24
$ws.Range("O24").Value = 23
$ws.Range("P24").Value = 100
$ws.Range("Q24").Value = 3
$ws.Range("R24").Value = 12
$ws.Range("S24").Value = 1
$ws.Range("T24").Value = 1
$ws.Range("W24").Value = 1
$ws.Range("X24").Value = 100
$ws.Range("Z24").Value = 0.1
$ws.Range("AA24").Value = 6
$ws.Range("AB24").Value = 200
$ws.Range("AE24").Value = 3000000

# Row 25
$ws.Range("O25").Value = 24
$ws.Range("P25").Value = 100
$ws.Range("Q25").Value = 3
$ws.Range("R25").Value = 12
$ws.Range("S25").Value = 1
$ws.Range("T25").Value = 1
$ws.Range("W25").Value = 1
$ws.Range("X25").Value = 100
$ws.Range("Z25").Value = 0.1
$ws.Range("AA25").Value = 6
$ws.Range("AB25").Value = 200
$ws.Range("AE25").Value = 3000000

# Row 26
$ws.Range("O26").Value = 25
$ws.Range("P26").Value = 100
$ws.Range("Q26").Value = 3
$ws.Range("R26").Value = 12
$ws.Range("S26").Value = 1
$ws.Range("T26").Value = 1
$ws.Range("W26").Value = 1
$ws.Range("X26").Value = 100
$ws.Range("Z26").Value = 0.1
$ws.Range("AA26").Value = 6
$ws.Range("AB26").Value = 200
$ws.Range("AE26").Value = 3000000

# Row 27
$ws.Range("O27").Value = 26
$ws.Range("P27").Value = 100
$ws.Range("Q27").Value = 3
$ws.Range("R27").Value = 12
$ws.Range("S27").Value = 1
$ws.Range("T27").Value = 1
$ws.Range("W27").Value = 1
$ws.Range("X27").Value = 100
$ws.Range("Z27").Value = 0.1
$ws.Range("AA27").Value = 6
$ws.Range("AB27").Value = 200
$ws.Range("AE27").Value = 3000000

# Row 28
$ws.Range("O28").Value = 27
$ws.Range("P28").Value = 100
$ws.Range("Q28").Value = 3
$ws.Range("R28").Value = 12
$ws.Range("S28").Value = 1
$ws.Range("T28").Value = 1
$ws.Range("W28").Value = 1
$ws.Range("X28").Value = 100
$ws.Range("Z28").Value = 0.1
$ws.Range("AA28").Value = 6
$ws.Range("AB28").Value = 200
$ws.Range("AE28").Value = 3000000

# Row 29
$ws.Range("O29").Value = 28
$ws.Range("P29").Value = 100
$ws.Range("Q29").Value = 3
$ws.Range("R29").Value = 12
$ws.Range("S29").Value = 1
$ws.Range("T29").Value = 1
$ws.Range("W29").Value = 1
$ws.Range("X29").Value = 100
$ws.Range("Z29").Value = 0.1
$ws.Range("AA29").Value = 6
$ws.Range("AB29").Value = 200
$ws.Range("AE29").Value = 3000000

# Row 30
$ws.Range("O30").Value = 29
$ws.Range("P30").Value = 100
$ws.Range("Q30").Value = 3
$ws.Range("R30").Value = 12
$ws.Range("S30").Value = 1
$ws.Range("T30").Value = 1
$ws.Range("W30").Value = 1
$ws.Range("X30").Value = 100
$ws.Range("Z30").Value = 0.1
$ws.Range("AA30").Value = 6
$ws.Range("AB30").Value = 200
$ws.Range("AE30").Value = 3000000

# Row 31
$ws.Range("O31").Value = 30
$ws.Range("P31").Value = 100
$ws.Range("Q31").Value = 3
$ws.Range("R31").Value = 12
$ws.Range("S31").Value = 1
$ws.Range("T31").Value = 1
$ws.Range("W31").Value = 1
$ws.Range("X31").Value = 100
$ws.Range("Z31").Value = 0.1
$ws.Range("AA31").Value = 6
$ws.Range("AB31").Value = 200
$ws.Range("AE31").Value = 3000000

# Row 32
$ws.Range("O32").Value = 31
$ws.Range("P32").Value = 100
$ws.Range("Q32").Value = 3
$ws.Range("R32").Value = 12
$ws.Range("S32").Value = 1
$ws.Range("T32").Value = 1
$ws.Range("W32").Value = 1
$ws.Range("X32").Value = 100
$ws.Range("Z32").Value = 0.1
$ws.Range("AA32").Value = 6
$ws.Range("AB32").Value = 200
$ws.Range("AE32").Value = 3000000

# Row 33
$ws.Range("O33").Value = 32
$ws.Range("P33").Value = 100
$ws.Range("Q33").Value = 3
$ws.Range("R33").Value = 12
$ws.Range("S33").Value = 1
$ws.Range("T33").Value = 1
$ws.Range("W33").Value = 1
$ws.Range("X33").Value = 100
$ws.Range("Z33").Value = 0.1
$ws.Range("AA33").Value = 6
$ws.Range("AB33").Value = 200
$ws.Range("AE33").Value = 3000000

# Update page setup to portrait orientation
$ws.PageSetup.Orientation = 1

# Update view: scroll so column O is the leftmost visible column, then select AB34
$excel.ActiveWindow.ScrollColumn = 15
$excel.ActiveWindow.ScrollRow = 1
$null = $ws.Range("AB34").Select()
